$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets the refreshed search values (2023-2024, semester 2, new end date)
$ws.Range("B2").Value = "2023-2024"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2"
$ws.Range("C2").ClearFormats()

$ws.Range("E2").Value = "31-03-2024"

# The old extra result rows are no longer part of the (narrowed) search results
$ws.Rows("3:4").Delete()
